$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "28.435.40", "1.000").
# Force Text format before assigning so Excel does not reinterpret them as
# numbers (which would drop thousands-grouping dots / trailing zeros / precision).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.435.40"
$ws.Range("E2").Value = "  +5.43%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.812.20"
$ws.Range("E3").Value = "  +4.07%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.12"
$ws.Range("E5").Value = "  +1.96%  "
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5512"
$ws.Range("E7").Value = "  +10.52%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3866"
$ws.Range("E8").Value = "  +7.98%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07615"
$ws.Range("E9").Value = "  +4.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.01"
$ws.Range("E10").Value = "  +1.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.133"
$ws.Range("E11").Value = "  +6.76%  "
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("E13").Value = "  +5.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.232"
$ws.Range("E14").Value = "  +4.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.338"
$ws.Range("E15").Value = "  +7.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.805.42"
$ws.Range("E16").Value = "  +4.06%  "
$ws.Range("E17").Value = "  +5.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001076"
$ws.Range("E18").Value = "  +4.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06477"
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.30"
$ws.Range("E21").Value = "  +4.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.000"
$ws.Range("E22").Value = "  +4.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.447.00"
$ws.Range("E23").Value = "  +5.27%  "
$ws.Range("E24").Value = "  +0.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.109"
$ws.Range("E25").Value = "  +2.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.75"
$ws.Range("E26").Value = "  +4.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.449"
$ws.Range("E27").Value = "  +14.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "156.53"
$ws.Range("E28").Value = "  +2.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.014.67"
$ws.Range("E29").Value = "  +4.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.22"
$ws.Range("E30").Value = "  +2.72%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.177"
$ws.Range("E31").Value = "  +10.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1036"
$ws.Range("E32").Value = "  +9.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.771"
$ws.Range("E33").Value = "  +7.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.642"
$ws.Range("E34").Value = "  +1.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.2317"
$ws.Range("E35").Value = "  +15.97%  "
$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.936"
$ws.Range("E36").Value = "  +19.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02334"
$ws.Range("E37").Value = "  +6.00%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06287"
$ws.Range("E38").Value = "  +6.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.64"
$ws.Range("E39").Value = "  +5.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6410"
$ws.Range("E40").Value = "  +6.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.048"
$ws.Range("E41").Value = "  +5.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.173"
$ws.Range("E42").Value = "  +5.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9994"
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.387"
$ws.Range("E44").Value = "  -2.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.51"
$ws.Range("E45").Value = "  +5.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6014"
$ws.Range("E46").Value = "  +6.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.693"
$ws.Range("E47").Value = "  +3.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.05"
$ws.Range("E48").Value = "  +3.32%  "
$ws.Range("E49").Value = "  +6.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.147"
$ws.Range("E50").Value = "  +4.20%  "
$ws.Range("E51").Value = "  +4.06%  "
